# Update "想去人数" (want-to-go count) values in column F for the
# "展览" and "全部类型" sheets, which carry identical data tables.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row -> new value for column F
$updates = @{
    4  = 2170
    14 = 13988
    26 = 637
    27 = 5106
    28 = 4
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
